$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("30:30").Insert()
$h = $ws.Rows("30:30").RowHeight
Write-Host "RowHeight after insert:" $h
$h29 = $ws.Rows("29:29").RowHeight
Write-Host "RowHeight 29:" $h29
